# Update sensitivity values for EOL-RIR Lifetime Min workbook
# Sheets: Neodymium, Dysprosium, Copper, Raw silicon
# Only Neodymium, Copper, and Raw silicon sheets change (B2:E5 each); Dysprosium unchanged.
# Note: literal scientific-notation numbers (e.g. 1.2E-05) are not parsed by this
# PowerShell engine, so values are supplied as strings and cast to [double].

$wb = $excel.ActiveWorkbook

# --- Neodymium sheet ---
$ws = $wb.Worksheets.Item("Neodymium")
$ws.Range("B2").Value = [double]"1.45568428466758E-19"
$ws.Range("C2").Value = [double]"3.234355465845293E-05"
$ws.Range("D2").Value = [double]"0.007019255690660976"
$ws.Range("E2").Value = [double]"0.009034683182445062"

$ws.Range("B3").Value = [double]"1.573017879546552E-09"
$ws.Range("C3").Value = [double]"0.000548952861088679"
$ws.Range("D3").Value = [double]"0.006513003445992318"
$ws.Range("E3").Value = [double]"0.0083886921140873"

$ws.Range("B4").Value = [double]"2.455525622174766E-11"
$ws.Range("C4").Value = [double]"0.0005091945878934708"
$ws.Range("D4").Value = [double]"0.005322594535634154"
$ws.Range("E4").Value = [double]"0.006745337844781117"

$ws.Range("C5").Value = [double]"1.125284018505872E-06"
$ws.Range("D5").Value = [double]"0.0003432694848831427"
$ws.Range("E5").Value = [double]"0.0005001255831081257"

# --- Copper sheet ---
$ws = $wb.Worksheets.Item("Copper")
$ws.Range("B2").Value = [double]"0.0004996343988088803"
$ws.Range("C2").Value = [double]"0.04114705065597703"
$ws.Range("D2").Value = [double]"0.7920825717866349"
$ws.Range("E2").Value = [double]"1.035122838457668"

$ws.Range("B3").Value = [double]"0.003689272420938101"
$ws.Range("C3").Value = [double]"0.04774564432274016"
$ws.Range("D3").Value = [double]"0.5565514789046994"
$ws.Range("E3").Value = [double]"0.7760924535686274"

$ws.Range("B4").Value = [double]"0.009318027231660879"
$ws.Range("C4").Value = [double]"0.03670670824608033"
$ws.Range("D4").Value = [double]"0.6625270780929425"
$ws.Range("E4").Value = [double]"0.8940617937696104"

$ws.Range("B5").Value = [double]"0.003016841357489923"
$ws.Range("C5").Value = [double]"0.03455053535924802"
$ws.Range("D5").Value = [double]"0.6243607850761044"
$ws.Range("E5").Value = [double]"0.8789148469914819"

# --- Raw silicon sheet ---
$ws = $wb.Worksheets.Item("Raw silicon")
$ws.Range("B2").Value = [double]"8.341175675907166E-05"
$ws.Range("C2").Value = [double]"0.001101398272791471"
$ws.Range("D2").Value = [double]"0.02423778873971561"
$ws.Range("E2").Value = [double]"0.02836977175650994"

$ws.Range("B3").Value = [double]"9.666419205255739E-05"
$ws.Range("C3").Value = [double]"0.001041167508083666"
$ws.Range("D3").Value = [double]"0.01323147715383524"
$ws.Range("E3").Value = [double]"0.01706756239850931"

$ws.Range("B4").Value = [double]"0.0005273422179231309"
$ws.Range("C4").Value = [double]"0.0009318616779647269"
$ws.Range("D4").Value = [double]"0.0178872485157633"
$ws.Range("E4").Value = [double]"0.02365167495720378"

$ws.Range("B5").Value = [double]"0.0002918342513006147"
$ws.Range("C5").Value = [double]"0.00107866079942882"
$ws.Range("D5").Value = [double]"0.02136331103838482"
$ws.Range("E5").Value = [double]"0.02530940646039109"
